{"js": "// Change the font size of the \"Val av testverktyget\" heading paragraph\n// from 20pt (w:sz/w:szCs = 40) to 14pt (w:sz/w:szCs = 28).\n//\n// The high-level Office.js Word.Font object only exposes `size`\n// (-> w:sz, the \"Western\" run size) and has no supported property for the\n// complex-script size (w:szCs). To flip both halves of the pair - exactly\n// like the OOXML diff requires - we read the paragraph's own OOXML,\n// surgically rewrite only the <w:sz>/<w:szCs> values that currently read\n// \"40\" to \"28\", and write the fragment back with insertOoxml(replace).\n// This keeps every other attribute (rsids, paraId, run props, etc.)\n// untouched.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst TARGET_TEXT = \"Val av testverktyget\";\nconst OLD_HALF_POINTS = \"40\"; // 20pt\nconst NEW_HALF_POINTS = \"28\"; // 14pt\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === TARGET_TEXT) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Paragraph \"' + TARGET_TEXT + '\" not found.');\n}\n\nconst ooxmlResult = target.getOoxml();\nawait context.sync();\n\n// Only touch <w:sz w:val=\"40\"/> and <w:szCs w:val=\"40\"/> (not any other\n// w:val=\"40\" that might belong to an unrelated element).\nconst sizeTagPattern = new RegExp(\n  '(<w:(?:sz|szCs)\\\\b[^>]*\\\\bw:val=\")' + OLD_HALF_POINTS + '(\")',\n  \"g\"\n);\nconst updatedXml = ooxmlResult.value.replace(\n  sizeTagPattern,\n  \"$1\" + NEW_HALF_POINTS + \"$2\"\n);\n\ntarget.insertOoxml(updatedXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Change the font size of the \"Val av testverktyget\" heading paragraph\n# from 20pt (w:sz/w:szCs = 40) to 14pt (w:sz/w:szCs = 28).\n#\n# Font.Size maps to the \"Western\" run size (w:sz); Font.SizeBi maps to the\n# complex-script size (w:szCs). Both need to be set to mirror the diff,\n# which changes both <w:sz> and <w:szCs> (in the paragraph-mark run\n# properties AND the text run's properties) from 40 to 28.\n\n$doc = $word.ActiveDocument\n\n$targetText = \"Val av testverktyget\"\n$newSizePoints = 14   # 28 half-points\n\n$rng = $doc.Content\n$rng.Find.ClearFormatting()\n\nwhile ($rng.Find.Execute($targetText) -and $rng.Find.Found) {\n    # Grab the whole paragraph (not just the matched text) so the\n    # paragraph-mark run properties get updated too, exactly like the\n    # source edit did for both the pPr/rPr and the r/rPr.\n    $para = $rng.Paragraphs(1)\n    $paraRange = $para.Range\n\n    $paraRange.Font.Size = $newSizePoints\n    $paraRange.Font.SizeBi = $newSizePoints\n\n    # Move past this paragraph before searching again, so we don't loop\n    # forever re-matching the same text.\n    $rng.SetRange($para.Range.End, $doc.Content.End)\n}\n"}
